$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.572.82"
Set-TextValue $ws.Range("E2") "  +2.32%  "
Set-TextValue $ws.Range("D3") "1.678.77"
Set-TextValue $ws.Range("E3") "  +2.88%  "
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "220.09"
Set-TextValue $ws.Range("E5") "  +2.63%  "
Set-TextValue $ws.Range("D6") "0.528"
Set-TextValue $ws.Range("E6") "  +2.06%  "
Set-TextValue $ws.Range("E7") "  -0.08%  "
Set-TextValue $ws.Range("D8") "30.04"
Set-TextValue $ws.Range("E8") "  +4.90%  "
Set-TextValue $ws.Range("E9") "  +2.43%  "
Set-TextValue $ws.Range("D10") "0.0630"
Set-TextValue $ws.Range("E10") "  +3.66%  "
Set-TextValue $ws.Range("E11") "  -0.76%  "
Set-TextValue $ws.Range("D12") "1.920.78"
Set-TextValue $ws.Range("E12") "  +2.97%  "
Set-TextValue $ws.Range("D13") "10.44"
Set-TextValue $ws.Range("E13") "  +12.76%  "
Set-TextValue $ws.Range("E14") "  +9.75%  "
Set-TextValue $ws.Range("D15") "1.683.06"
Set-TextValue $ws.Range("E15") "  +3.17%  "
Set-TextValue $ws.Range("D16") "3.98"
Set-TextValue $ws.Range("E16") "  +3.03%  "
Set-TextValue $ws.Range("D17") "30.571.00"
Set-TextValue $ws.Range("E17") "  +2.30%  "
Set-TextValue $ws.Range("D18") "66.41"
Set-TextValue $ws.Range("E18") "  +3.37%  "
Set-TextValue $ws.Range("D19") "245.21"
Set-TextValue $ws.Range("E19") "  +1.15%  "
Set-TextValue $ws.Range("D20") "0.0₃0717"
Set-TextValue $ws.Range("E20") "  +2.43%  "
Set-TextValue $ws.Range("E21") "  -0.10%  "
Set-TextValue $ws.Range("E24") "  +1.09%  "
Set-TextValue $ws.Range("D25") "157.76"
Set-TextValue $ws.Range("E25") "  +0.17%  "
Set-TextValue $ws.Range("D26") "15.90"
Set-TextValue $ws.Range("E26") "  +2.43%  "
Set-TextValue $ws.Range("E27") "  +2.20%  "
Set-TextValue $ws.Range("E28") "  +1.87%  "
Set-TextValue $ws.Range("E29") "  -0.11%  "
Set-TextValue $ws.Range("D30") "0.0497"
Set-TextValue $ws.Range("E30") "  +2.44%  "
Set-TextValue $ws.Range("E31") "  +2.65%  "
Set-TextValue $ws.Range("E32") "  +3.48%  "
Set-TextValue $ws.Range("D33") "1.509.50"
Set-TextValue $ws.Range("E33") "  +5.58%  "
Set-TextValue $ws.Range("E34") "  +4.25%  "
Set-TextValue $ws.Range("D35") "1.76"
Set-TextValue $ws.Range("E35") "  +7.24%  "
Set-TextValue $ws.Range("E36") "  -0.12%  "
Set-TextValue $ws.Range("D37") "83.72"
Set-TextValue $ws.Range("E37") "  +11.05%  "
Set-TextValue $ws.Range("E38") "  +5.54%  "
Set-TextValue $ws.Range("D39") "0.592"
Set-TextValue $ws.Range("E39") "  +7.66%  "
Set-TextValue $ws.Range("D40") "2.71"
Set-TextValue $ws.Range("E40") "  -2.92%  "
Set-TextValue $ws.Range("E41") "  -0.12%  "
Set-TextValue $ws.Range("D42") "0.839"
Set-TextValue $ws.Range("E42") "  +1.75%  "
Set-TextValue $ws.Range("E43") "  +2.10%  "
Set-TextValue $ws.Range("E44") "  -0.20%  "
Set-TextValue $ws.Range("E45") "  +0.52%  "
Set-TextValue $ws.Range("E46") "  -0.05%  "
Set-TextValue $ws.Range("D49") "1.813.90"
Set-TextValue $ws.Range("E49") "  +2.24%  "
Set-TextValue $ws.Range("D50") "94.70"
Set-TextValue $ws.Range("E50") "  +6.34%  "
Set-TextValue $ws.Range("E51") "  +2.14%  "

# Row 22/23 swap: Avalanche <-> Uniswap
Set-TextValue $ws.Range("B22") "Uniswap"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "4.28"
Set-TextValue $ws.Range("E22") "  +3.83%  "
Set-TextValue $ws.Range("B23") "Avalanche"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D23") "10.15"
Set-TextValue $ws.Range("E23") "  +3.18%  "

# Row 47/48 swap: FraxShare <-> BitcoinSV
Set-TextValue $ws.Range("B47") "BitcoinSV"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D47") "51.68"
Set-TextValue $ws.Range("E47") "  -3.56%  "
Set-TextValue $ws.Range("B48") "FraxShare"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D48") "5.55"
Set-TextValue $ws.Range("E48") "  +3.18%  "
